$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching original inline-string formatting)
$textCells = @("D5", "D6", "D7", "D9", "D10", "D13", "D14", "D16", "D17", "D20", "D23", "D24", "D28", "D29", "D31", "D33", "D35", "D39", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "36.396.89"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.014.31"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "252.17"
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("D7").Value = "62.32"
$ws.Range("E7").Value = "  +10.14%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.371"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "58.43"
$ws.Range("E10").Value = "  -6.44%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "0.905"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "14.94"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "2.307.84"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "20.74"
$ws.Range("E16").Value = "  +16.50%  "
$ws.Range("D17").Value = "5.48"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "2.021.44"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "36.340.54"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "72.07"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "234.67"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  +19.29%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "163.52"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "19.65"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "5.13"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  +20.86%  "
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("D35").Value = "0.0609"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "5.92"
$ws.Range("E39").Value = "  +16.14%  "
$ws.Range("E40").Value = "  +16.11%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "17.10"
$ws.Range("E43").Value = "  +7.20%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0216"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.00"
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  +17.58%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.456.12"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "95.27"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "47.26"
$ws.Range("E51").Value = "  +1.96%  "
